# Applies the zh-Hant -> zh-Hans (Simplified Chinese) re-translation
# described by the commit "New translations email 5-1 [template] partner
# email - invite revoked.docx (Chinese Simplified)".

$d = $word.ActiveDocument

function Replace-All($find, $replace) {
    $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $replace, 2) | Out-Null
}

# Hyperlink text "英語" / trailing language list, appears twice with the
# same target replacement text, so a global replace covers both.
Replace-All "英語" "英语"
Replace-All " / 葡萄牙語 / 法語 / 泰語 / 越南語 / 西班牙語" " / 葡萄牙语 / 法语 / 泰语 / 越南语 / 西班牙语"

# Table: summary / target-audience block
Replace-All "簡介" "简要"
Replace-All "發送給在目標國家的合作夥伴的電子郵件，這些合作夥伴已回應參加，但在截止日期前未提交文件。 我們將取消他們的邀請。 將通過 customer.io 發送" "发给在目标国家已确认出席但未在截止日期前提交文件的合作伙伴的邮件。 我们将撤回他们的邀请。 将通过 customer.io 发送"
Replace-All "目標受眾" "目标受众"
Replace-All "未按時提交文件的被邀請合作夥伴" "未及时提交文件的被邀请合作伙伴"

# Subject line
Replace-All "主題行" "主题行"
Replace-All "[活動名稱]" "[活动名称]"
Replace-All " 註冊" " 注册"

# Heading
Replace-All "沒有及時收到您的文件" "没有及时收到您的文件"

# Greeting / body — "[合作夥伴姓名]" becomes the untranslated English
# placeholder "[PARTNER NAME]"
Replace-All "[合作夥伴姓名]" "[PARTNER NAME]"

Replace-All "截止日期（" "We didn’t receive your documents by the deadline ("
Replace-All "[日月年]" "[DD Mmm YYYY]"
Replace-All "）前沒有收到您的文件。 很遺憾，無法為您辦理 " "). 很遗憾，无法为您办理 "
Replace-All " 的註冊手續。" " 的注册手续。"

Replace-All "衷心祝愿您一切順利，並希望在下一次 " "衷心祝愿您一切顺利，并希望在下一次 "
Replace-All "會議/研討會/聯盟會員旅行" "会议/研讨会/联盟会员旅行"
Replace-All "中見到您。" "中见到您。"

Replace-All "如有任何疑問，請通過 " "如有任何疑问，请通过 "
Replace-All "[電子郵件地址]" "[电子邮件地址]"
Replace-All "[WHATSAPP 號碼]" "[WHATSAPP 号码]"
Replace-All " (WhatsApp) 聯繫您的區域經理 " " (WhatsApp) 联系您的区域经理 "

# Comments (not reachable through $d.Content — update via the Comments
# collection directly).
for ($i = 1; $i -le $d.Comments.Count; $i++) {
    $c = $d.Comments.Item($i)
    if ($c.Range.Text -eq "選擇其中一個") {
        $c.Range.Text = "选择任一"
    }
}
